$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - requirement #4
$ws.Range("B6").NumberFormat = "General"
$ws.Range("B6").Value = 4
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "req4"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "Software requirement description four"

# Row 7 - requirement #5
$ws.Range("B7").NumberFormat = "General"
$ws.Range("B7").Value = 5
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "req5"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "Software requirement description five"
